$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.374.73'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '1.566.41'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9981'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3776'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.86%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.40'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3407'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07608'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.140'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("E13").Value = '  -0.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.980'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.937'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001132'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.552.71'
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06734'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.204'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("D24").Value = '22.362.77'
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.403'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.709'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.028'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.27'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("D31").Value = '1.748.12'
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.016'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.108'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9899'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '10.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08472'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.403'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02515'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2297'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06479'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.400'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.35%  '
$ws.Range("E42").Value = '  -2.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6317'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.97'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.812'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.69%  '
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.083'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("E49").Value = '  +0.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '124.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("E51").Value = '  +0.48%  '
